$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Semana_14")

# Row 5
$ws.Range("L5").Value = 5
$ws.Range("M5").Value = 1

# Row 11
$ws.Range("L11").Value = 2
$ws.Range("M11").Value = 1

# Row 14
$ws.Range("L14").Value = 12
$ws.Range("M14").Value = -2

# Row 46
$ws.Range("L46").Value = 10
$ws.Range("M46").Value = 2

# Row 50
$ws.Range("L50").Value = 13
$ws.Range("M50").Value = 1

# Row 53
$ws.Range("L53").Value = 5
$ws.Range("M53").Value = 1

# Row 56
$ws.Range("L56").Value = 2
$ws.Range("M56").Value = 1

# Row 69
$ws.Range("L69").Value = 5
$ws.Range("M69").Value = 1

# Summary metrics
$ws.Range("C82").Value = 243
$ws.Range("C93").Value = 6
